# Resident Services_Requirements.xlsx - add "Reg Proc" column (T) to the
# "Details" sheet table (Table2) and populate new clarification notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$lo = $ws.ListObjects.Item("Table2")

# Add the new table column (extends Table2 from A2:S16 to A2:T16, updates
# the AutoFilter range and the worksheet dimension automatically).
$newCol = $lo.ListColumns.Add()
$newCol.Name = "Reg Proc"

# --- Header cell (row 2) ---------------------------------------------
$ws.Range("S2").Copy() | Out-Null
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$ws.Range("T2").Value = "Reg Proc"

# --- Row 4 -------------------------------------------------------------
$ws.Range("S5").Copy() | Out-Null
$ws.Range("T4").PasteSpecial(-4122) | Out-Null
$ws.Range("T4").Value = "When UIN IS needed to be generated" + [char]10 + "1.the Acknowledgment from Print queue- what needs to be done" + [char]10 + "Time period " + [char]10 + "2. If there is a print failure- no need to handle from MOSIP" + [char]10 + "User Story ?"

# --- Row 5 -------------------------------------------------------------
$ws.Range("S5").Copy() | Out-Null
$ws.Range("T5").PasteSpecial(-4122) | Out-Null
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor" + [char]10 + "ID Repo- Might not be there in ID Repo as well"

# --- Row 6 -------------------------------------------------------------
$ws.Range("S4").Copy() | Out-Null
$ws.Range("T6").PasteSpecial(-4122) | Out-Null
$ws.Range("T6").Value = "ID Repo- need to know "

# --- Row 7 -------------------------------------------------------------
$ws.Range("S5").Copy() | Out-Null
$ws.Range("T7").PasteSpecial(-4122) | Out-Null
$ws.Range("T7").Value = "there shud be a label as Res_Service" + [char]10 + "Reg Client packet needs to be understood" + [char]10 + "Service from Reg proc needs to be developed"

# --- Row 8 (existing S8 text is replaced, new T8 is added) ------------
$ws.Range("S5").Copy() | Out-Null
$ws.Range("S8").PasteSpecial(-4122) | Out-Null
$ws.Range("S8").Value = "Reg proc" + [char]10 + "Archival policy"

$ws.Range("S5").Copy() | Out-Null
$ws.Range("T8").PasteSpecial(-4122) | Out-Null
$ws.Range("T8").Value = "Under processing" + [char]10 + "Processed" + [char]10

# --- Row 9 -------------------------------------------------------------
$ws.Range("S5").Copy() | Out-Null
$ws.Range("T9").PasteSpecial(-4122) | Out-Null
$ws.Range("T9").Value = "Under processing" + [char]10 + "Processed"

# --- Row 10 --------------------------------------------------------------
$ws.Range("S4").Copy() | Out-Null
$ws.Range("T10").PasteSpecial(-4122) | Out-Null
$ws.Range("T10").Value = "E-UIN Generation"

# --- Column width for the new column T (bestFit/autofit width) ---------
$ws.Columns.Item(20).ColumnWidth = 31.25

# --- View: refresh the active cell / selection on this sheet ------------
$ws.Activate()
$ws.Range("T4").Select()

$excel.CutCopyMode = $false
